$d = $word.ActiveDocument

# Change 1: merge "SENDPLAYERLIST" + "#No Players" runs (removing _GoBack bookmark between them)
$d.Content.Find.Execute("SENDPLAYERLIST#No Players", $true, $false, $false, $false, $false, $true, 1, $false, "SENDPLAYERLIST#No Players", 2)

# Change 2: JUMPTHISTURN -> QUITTHISMATCH
$d.Content.Find.Execute("JUMPTHISTURN", $true, $false, $false, $false, $false, $true, 1, $false, "QUITTHISMATCH", 2)

# Change 2b: text
$d.Content.Find.Execute("Der Client teilt dem Server mit, dass er einmal aussetzen möchte.", $true, $false, $false, $false, $false, $true, 1, $false, "Der Client teilt dem Server mit, dass er diesen Match beenden möchte.", 2)
